$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.952.02"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.881.23"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'587.14"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -5.74%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -5.08%  "
$ws.Range("D11").Value = "'0.426"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("D12").Value = "'0.0000216"
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("D13").Value = "'32.13"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "3.357.73"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "60.837.10"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "2.879.37"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "'6.46"
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("D19").Value = "'423.73"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "'13.20"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").Value = "'0.650"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("D22").Value = "'6.90"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").Value = "'79.65"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "'10.31"
$ws.Range("E24").Value = "  -5.35%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -7.60%  "
$ws.Range("D27").Value = "'11.32"
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").Value = "'2.05"
$ws.Range("E29").Value = "  -9.50%  "
$ws.Range("D30").Value = "'6.59"
$ws.Range("E30").Value = "  -6.08%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'25.50"
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("E33").Value = "  -6.02%  "
$ws.Range("D34").Value = "0.0₃0837"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "'0.966"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("E36").Value = "  -4.44%  "
$ws.Range("D37").Value = "'48.97"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'2.77"
$ws.Range("E38").Value = "  -7.91%  "
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").Value = "'8.30"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("D41").Value = "'0.114"
$ws.Range("E41").Value = "  -6.66%  "
$ws.Range("E42").Value = "  -7.40%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.660.65"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "'37.65"
$ws.Range("E44").Value = "  -7.32%  "
$ws.Range("D45").Value = "'131.67"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "'0.0327"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("D47").Value = "'346.37"
$ws.Range("E47").Value = "  -8.91%  "
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("D50").Value = "'22.13"
$ws.Range("E50").Value = "  -7.11%  "
$ws.Range("E51").Value = "  -5.05%  "
